$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 55 and 58 had their observation-specific data swapped (the two
# occurrence records traded places while a handful of shared/location columns
# - P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY - stayed
# put). Apply the resulting target values directly.
# ---------------------------------------------------------------------------

# Row 55 (becomes the former row 58 occurrence: Rynkskinn / Phlebia centrifuga)
$ws.Range("A55").Value = 111898507
$ws.Range("B55").Value = 89845
$ws.Range("D55").Value = "VU"
$ws.Range("E55").Value = 1209
$ws.Range("F55").Value = "Rynkskinn"
$ws.Range("G55").Value = "Phlebia centrifuga"
$ws.Range("H55").Value = "P.Karst."
$ws.Range("K55").Value = ""
$ws.Range("L55").Value = ""
$ws.Range("Q55").Value = 650086.8716060545
$ws.Range("R55").Value = 6654015.064976334
$ws.Range("AH55").Value = "Ängsblandskog"
$ws.Range("AI55").Value = ""
$ws.Range("AJ55").Value = "gran"
$ws.Range("AK55").Value = "Picea abies"
$ws.Range("AM55").Value = "Liggande död trädstam, utan markontakt"
$ws.Range("AO55").Value = "Horizontal, dead without ground contact # Picea abies"

# Row 58 (becomes the former row 55 occurrence: Blåsippa / Hepatica nobilis)
$ws.Range("A58").Value = 111898889
$ws.Range("B58").Value = 98535
$ws.Range("D58").Value = "LC"
$ws.Range("E58").Value = 222498
$ws.Range("F58").Value = "Blåsippa"
$ws.Range("G58").Value = "Hepatica nobilis"
$ws.Range("H58").Value = "Schreb."
$ws.Range("K58").Value = "fullt utvecklade blad"
$ws.Range("L58").Value = ""
$ws.Range("Q58").Value = 650135.0421630922
$ws.Range("R58").Value = 6654002.501842719
$ws.Range("AH58").Value = "Ängsbarrskog"
$ws.Range("AI58").Value = "Ungskog"
$ws.Range("AJ58").Value = ""
$ws.Range("AK58").Value = ""
$ws.Range("AM58").Value = ""
$ws.Range("AO58").Value = ""

# ---------------------------------------------------------------------------
# Rows 60 and 61 (both Knärot / Goodyera repens) swapped their Id, Antal and
# coordinate values.
# ---------------------------------------------------------------------------

$ws.Range("A60").Value = 111911698
$ws.Range("I60").NumberFormat = "@"
$ws.Range("I60").Value = "16"
$ws.Range("Q60").Value = 650032.9755174413
$ws.Range("R60").Value = 6654279.303373625

$ws.Range("A61").Value = 111911660
$ws.Range("I61").NumberFormat = "@"
$ws.Range("I61").Value = "19"
$ws.Range("Q61").Value = 650026.652882754
$ws.Range("R61").Value = 6654299.07778531
